$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing D-column figures (weekly counts refreshed from source) ---
$ws.Range("D2").Value = 11756
$ws.Range("D3").Value = 11501
$ws.Range("D7").Value = 11860
$ws.Range("D8").Value = 11148
$ws.Range("D12").Value = 11968
$ws.Range("D13").Value = 10646
$ws.Range("D17").Value = 12012
$ws.Range("D18").Value = 10095
$ws.Range("D22").Value = 12069
$ws.Range("D23").Value = 9666
$ws.Range("D24").Value = 2403
$ws.Range("D26").Value = 1883
$ws.Range("D27").Value = 12111
$ws.Range("D28").Value = 9207
$ws.Range("D29").Value = 2904
$ws.Range("D31").Value = 2295
$ws.Range("D32").Value = 12151
$ws.Range("D33").Value = 8799
$ws.Range("D34").Value = 3352
$ws.Range("D36").Value = 2646
$ws.Range("D37").Value = 12189
$ws.Range("D38").Value = 8275
$ws.Range("D39").Value = 3914
$ws.Range("D41").Value = 3113
$ws.Range("D42").Value = 12224
$ws.Range("D43").Value = 7708
$ws.Range("D44").Value = 4516
$ws.Range("D46").Value = 3614
$ws.Range("D47").Value = 12249
$ws.Range("D48").Value = 7160
$ws.Range("D49").Value = 5089
$ws.Range("D51").Value = 4074
$ws.Range("D52").Value = 12277
$ws.Range("D53").Value = 6518
$ws.Range("D54").Value = 5759
$ws.Range("D56").Value = 4578
$ws.Range("D57").Value = 12300
$ws.Range("D58").Value = 5908
$ws.Range("D59").Value = 6392
$ws.Range("D61").Value = 5104
$ws.Range("D62").Value = 12311
$ws.Range("D63").Value = 5742
$ws.Range("D64").Value = 6569
$ws.Range("D66").Value = 5252
$ws.Range("D67").Value = 12333
$ws.Range("D68").Value = 5558
$ws.Range("D69").Value = 6775
$ws.Range("D71").Value = 5410
$ws.Range("D72").Value = 12349
$ws.Range("D73").Value = 5122
$ws.Range("D74").Value = 7227
$ws.Range("D76").Value = 5811
$ws.Range("D77").Value = 12371
$ws.Range("D78").Value = 4644
$ws.Range("D79").Value = 7727
$ws.Range("D81").Value = 6255
$ws.Range("D82").Value = 12389
$ws.Range("D83").Value = 4125
$ws.Range("D84").Value = 8264
$ws.Range("D86").Value = 6756
$ws.Range("D87").Value = 12416
$ws.Range("D88").Value = 3719
$ws.Range("D89").Value = 8697
$ws.Range("D91").Value = 7159
$ws.Range("D92").Value = 12432
$ws.Range("D93").Value = 3381
$ws.Range("D94").Value = 9051
$ws.Range("D96").Value = 7517
$ws.Range("D97").Value = 12453
$ws.Range("D98").Value = 3070
$ws.Range("D99").Value = 9383
$ws.Range("D101").Value = 7834
$ws.Range("D102").Value = 12468
$ws.Range("D103").Value = 2794
$ws.Range("D104").Value = 9674
$ws.Range("D105").Value = 1549
$ws.Range("D106").Value = 8125

# --- Remove now-unused trailing blank rows (112:131), shifting nothing below them ---
$ws.Range("A112:D131").EntireRow.Delete() | Out-Null

# --- Append the new ISO week 2025-09 block in rows 107:111 ---
$ws.Range("A107").Value = 202509
$ws.Range("B107").Value = 45718
$ws.Range("C107").Value = "farms_total_count"
$ws.Range("D107").Value = 12479

$ws.Range("A108").Value = 202509
$ws.Range("B108").Value = 45718
$ws.Range("C108").Value = "farms_to_examine_count"
$ws.Range("D108").Value = 2508

$ws.Range("A109").Value = 202509
$ws.Range("B109").Value = 45718
$ws.Range("C109").Value = "farms_examined_count"
$ws.Range("D109").Value = 9971

$ws.Range("A110").Value = 202509
$ws.Range("B110").Value = 45718
$ws.Range("C110").Value = "farms_examined_positive_count"
$ws.Range("D110").Value = 1541

$ws.Range("A111").Value = 202509
$ws.Range("B111").Value = 45718
$ws.Range("C111").Value = "farms_examined_negative_count"
$ws.Range("D111").Value = 8430

# --- Restore the view to the cursor position saved with the workbook ---
$ws.Range("C10").Select() | Out-Null
